$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style of the existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Data for columns I (I0) and J (IF), keyed by row number
$values = @{
    2  = @(1, 6)
    3  = @(1, 6)
    4  = @(6, 7)
    5  = @(1, 6)
    6  = @(1, 6)
    7  = @(1, 5)
    8  = @(1, 6)
    9  = @(1, 6)
    10 = @(1, 5)
    11 = @(6, 8)
    12 = @(1, 6)
    13 = @(1, 6)
    14 = @(1, 6)
    15 = @(1, 6)
    16 = @(1, 5)
    17 = @(1, 5)
    18 = @(1, 4)
    19 = @(1, 5)
    20 = @(1, 7)
    21 = @(1, 7)
    22 = @(1, 5)
    23 = @(1, 7)
    24 = @(1, 4)
    25 = @(1, 4)
    26 = @(1, 6)
    27 = @(1, 6)
    28 = @(1, 4)
    29 = @(6, 8)
    30 = @(1, 6)
    31 = @(1, 5)
    32 = @(3, 6)
    33 = @(5, 7)
    34 = @(6, 7)
    35 = @(6, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
